# Adds a new "2022-Q3" worksheet (with its fund holdings data) right after
# "总计" and before "2022-Q2", and updates the "总计" summary sheet with a
# new leading row for the 2022-Q3 totals (shifting the existing rows down).

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $addr, $text) {
    # Force a value to be stored as text even when it looks numeric
    # (e.g. "159617", "2.93", "0.0410"), mirroring how the source workbook
    # keeps these as literal strings instead of coerced numbers. Using a
    # leading apostrophe marks the input as text; resetting the style back
    # to "Normal" afterwards drops the quote-prefix formatting flag so the
    # cell ends up with the same (default) style as the rest of the sheet.
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet by duplicating "2022-Q2" (so it keeps
#    identical sheet-level formatting/structure), then rename it and move
#    it immediately before "2022-Q2".
# ---------------------------------------------------------------------
$existingQ2 = $wb.Worksheets.Item("2022-Q2")
$existingQ2.Copy($existingQ2, $null)

$newSheet = $wb.Worksheets.Item("2022-Q2 (2)")
$newSheet.Name = "2022-Q3"

# The source sheet had 3 data rows; the new sheet only needs 2, so drop
# the extra (4th) row entirely.
$newSheet.Rows.Item(4).Delete()

# ---------------------------------------------------------------------
# 2. Populate the new sheet's data (header row is already correct from
#    the copy).
# ---------------------------------------------------------------------
$newSheet.Range("A2").Value = 0
Set-TextCell $newSheet "B2" "159617"
Set-TextCell $newSheet "C2" "华夏中证智选500价值稳健策略ETF"
Set-TextCell $newSheet "D2" "2.93"
Set-TextCell $newSheet "E2" "97.05"
Set-TextCell $newSheet "F2" "1.40"
Set-TextCell $newSheet "G2" "0.0410"
$newSheet.Range("H2").Value = 8

$newSheet.Range("A3").Value = 1
Set-TextCell $newSheet "B3" "512590"
Set-TextCell $newSheet "C3" "浦银安盛中证高股息精选ETF"
Set-TextCell $newSheet "D3" "0.45"
Set-TextCell $newSheet "E3" "90.87"
Set-TextCell $newSheet "F3" "2.27"
Set-TextCell $newSheet "G3" "0.0102"
$newSheet.Range("H3").Value = 3

# ---------------------------------------------------------------------
# 3. Update the "总计" summary sheet: insert a new row below the header
#    for 2022-Q3 and shift the existing quarters' rows down by one.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

$summary.Range("A2").Value = 0
Set-TextCell $summary "B2" "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.05

# The newly-inserted row pushed every other row down by one, so column A's
# running index needs to be bumped back into sequence (0..7) and the last
# (previously non-existent) row needs to be filled in with the 2020-Q4
# figures that fell off the bottom.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
$summary.Range("A8").Value = 6

$summary.Range("A9").Value = 7
Set-TextCell $summary "B9" "2020-Q4"
$summary.Range("C9").Value = 9
$summary.Range("D9").Value = 0.18

# Make sure the summary sheet is the active one, matching the original
# workbook's default view.
$summary.Activate()
